$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-12-14 Saturday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-12-15 Sunday", 2) | Out-Null

# Update the division-problem table by addressing cells directly
# (several values repeat, so text-based Find/Replace would be ambiguous)
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "67÷6=11, 1"
$t.Cell(1, 2).Range.Text = "93÷7=13, 2"
$t.Cell(1, 3).Range.Text = "48÷8=6, 0"
$t.Cell(1, 4).Range.Text = "14÷3=4, 2"
$t.Cell(1, 5).Range.Text = "10÷8=1, 2"

# Row 5
$t.Cell(5, 1).Range.Text = "98÷5=19, 3"
$t.Cell(5, 2).Range.Text = "96÷3=32, 0"
$t.Cell(5, 3).Range.Text = "37÷2=18, 1"
$t.Cell(5, 4).Range.Text = "25÷8=3, 1"
$t.Cell(5, 5).Range.Text = "53÷7=7, 4"

# Row 9
$t.Cell(9, 1).Range.Text = "63÷9=7, 0"
$t.Cell(9, 2).Range.Text = "19÷9=2, 1"
$t.Cell(9, 3).Range.Text = "10÷3=3, 1"
$t.Cell(9, 4).Range.Text = "21÷3=7, 0"
$t.Cell(9, 5).Range.Text = "81÷7=11, 4"

# Row 13
$t.Cell(13, 1).Range.Text = "46÷8=5, 6"
$t.Cell(13, 2).Range.Text = "94÷9=10, 4"
$t.Cell(13, 3).Range.Text = "56÷9=6, 2"
$t.Cell(13, 4).Range.Text = "94÷4=23, 2"
$t.Cell(13, 5).Range.Text = "99÷3=33, 0"

# Row 17
$t.Cell(17, 1).Range.Text = "23÷5=4, 3"
$t.Cell(17, 2).Range.Text = "58÷5=11, 3"
$t.Cell(17, 3).Range.Text = "93÷9=10, 3"
$t.Cell(17, 4).Range.Text = "82÷5=16, 2"
$t.Cell(17, 5).Range.Text = "47÷5=9, 2"
